# xls export geometry fix, added server power\health status
#
# Rename a few report headers (shared strings) and shrink/retarget
# some column widths to better fit the new, shorter labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F: "System memory size" -> "Memory tot.size"
$ws.Range("F1").Value = "Memory tot.size"

# Column H: "Memory module part number" -> "Memory P/Ns"
$ws.Range("H1").Value = "Memory P/Ns"

# Column M: "HDD slot population" -> "HDD slot pop."
$ws.Range("M1").Value = "HDD slot pop."

# Column N: "PSU part number" -> "PSU P/Ns"
$ws.Range("N1").Value = "PSU P/Ns"

# Adjust column widths to match the new, narrower header text
# (values chosen so the saved width lands as close as possible to the
# target widths: F=15.71, H=11.71, M=13.71, N=8.71 "characters").
$ws.Columns.Item(6).ColumnWidth = 14.83
$ws.Columns.Item(8).ColumnWidth = 10.83
$ws.Columns.Item(13).ColumnWidth = 12.83
$ws.Columns.Item(14).ColumnWidth = 7.83
